$wb = $excel.ActiveWorkbook

# --- Sheet "Typography": set J4 (Widget Wildcard Characters) ---
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("J4").Value = "-., 0123456789"

# --- Sheet "Translation": fill rows 7-14 with new text entries ---
$wsTrans = $wb.Worksheets.Item("Translation")

# Match formatting of existing data rows (no explicit cell style override)
$wsTrans.Range("B7:F14").Style = "Normal"

$rows = @(
    @{ Row=7;  B="SingleUseId3";  C="Default"; D="Center"; E="<>"; F="LTR" },
    @{ Row=8;  B="SingleUseId4";  C="Default"; D="Left";   E="<>"; F="LTR" },
    @{ Row=9;  B="SingleUseId5";  C="Default"; D="Right";  E="<>"; F="LTR" },
    @{ Row=10; B="SingleUseId6";  C="Default"; D="Left";   E="<>"; F="LTR" },
    @{ Row=11; B="SingleUseId7";  C="Default"; D="Left";   E="<>"; F="LTR" },
    @{ Row=12; B="SingleUseId8";  C="Default"; D="Right";  E="<>"; F="LTR" },
    @{ Row=13; B="SingleUseId9";  C="Default"; D="Left";   E="<>"; F="LTR" },
    @{ Row=14; B="SingleUseId10"; C="Default"; D="Right";  E="<>"; F="LTR" }
)

foreach ($r in $rows) {
    $wsTrans.Cells.Item($r.Row, 2).Value = $r.B
    $wsTrans.Cells.Item($r.Row, 3).Value = $r.C
    $wsTrans.Cells.Item($r.Row, 4).Value = $r.D
    $wsTrans.Cells.Item($r.Row, 5).Value = $r.E
    $wsTrans.Cells.Item($r.Row, 6).Value = $r.F
}
